$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Päivitetty:" label and date + name
$ws.Range("F1").Value = "Päivitetty:"
$ws.Range("G1").NumberFormat = "mm-dd-yy"
$ws.Range("G1").Value = Get-Date -Year 2017 -Month 11 -Day 27 -Hour 0 -Minute 0 -Second 0
$ws.Range("H1").Value = "(Mikko)"
$ws.Range("H1").HorizontalAlignment = -4152

# Update F5 status from "IN PROGRESS" to "DONE"
$ws.Range("F5").Value = "DONE"

# Update selection
$ws.Range("H4").Select()
